# Applies the cryptos.xlsx price/volume/coin refresh described in the commit
# 'Updated cryptos list on Sun Feb 26 18:47:31 UTC 2023 with GitHub Actions'.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '23.373.23'
$ws.Range("E2").Value = '  +1.33%  '

# Row 3
$ws.Range("D3").Value = '1.617.90'
$ws.Range("E3").Value = '  +1.56%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9953'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.58%  '

# Row 5
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.05'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.34%  '

# Row 6
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9981'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.35%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3770'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.17%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '53.23'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +5.08%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3637'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.80%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.273'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.25%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08175'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.85%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9983'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.43%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.14'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +4.37%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.643'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.95%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.396'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.48%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001252'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.54%  '

# Row 17
$ws.Range("D17").Value = '1.608.15'
$ws.Range("E17").Value = '  +1.06%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.34'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.99%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06931'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.41%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.30'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.30%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.565'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.20%  '

# Row 22
$ws.Range("E22").Value = '  -0.06%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.95'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.05%  '

# Row 24
$ws.Range("D24").Value = '23.304.37'
$ws.Range("E24").Value = '  +1.07%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.108'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +10.52%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.409'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.41%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.34'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.56%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.67'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.30%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.286'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.95%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '135.60'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.37%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.408'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.34%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.844'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.57%  '

# Row 33
$ws.Range("D33").Value = '1.787.24'
$ws.Range("E33").Value = '  +1.33%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9636'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.13%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02774'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.25%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.38'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.65%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.07389'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.72%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.182'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.07%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2522'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.00%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08805'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.03%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.399'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.92%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7123'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.10%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.65'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +3.08%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.89'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +5.39%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6566'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.85%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.345'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.31%  '

# Row 47
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9977'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.27%  '

# Row 48
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.021'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.27%  '

# Row 49
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '132.74'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.82%  '

# Row 50
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07998'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.09%  '

# Row 51
$ws.Range("B51").Value = 'Flow'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.200'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.78%  '
